# création d'un alerte quand on créé un fichier de données
# Met à jour les lignes de présence existantes et ajoute les nouvelles
# lignes de suivi de présence dans la feuille.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force le stockage en texte pour éviter qu'Excel ne réinterprète les
    # dates (ex: "2021-02-01") comme des numéros de série de date, puis
    # retire le format appliqué pour ne laisser aucune trace de style.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Ligne 2 (ID 1) : Njohou Landry - cours de Maths, semestre 1, present
$ws.Range("B2").Value = "Njohou"
$ws.Range("C2").Value = "Landry"
$ws.Range("E2").Value = "Maths"

# Ligne 3 (ID 2) : Nkoa Christophe - cours d'Electronique, semestre 2, en retard
$ws.Range("B3").Value = "Nkoa"
$ws.Range("C3").Value = "Christophe"
Set-TextValue $ws.Range("D3") "2021-01-30"
$ws.Range("E3").Value = "Electronique"
$ws.Range("F3").Value = "semestre 2"
$ws.Range("G3").Value = 2019
$ws.Range("H3").Value = "En retard"

# Ligne 4 (ID 3) : Njohou Landry - cours d'Electronique, semestre 2, en retard
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Njohou"
$ws.Range("C4").Value = "Landry"
Set-TextValue $ws.Range("D4") "2021-02-04"
$ws.Range("E4").Value = "Electronique"
$ws.Range("F4").Value = "semestre 2"
$ws.Range("G4").Value = 2019
$ws.Range("H4").Value = "En retard"

# Ligne 5 (ID 4) : Njohou Landry - cours d'Electronique, semestre 2, absent
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Njohou"
$ws.Range("C5").Value = "Landry"
Set-TextValue $ws.Range("D5") "2021-01-14"
$ws.Range("E5").Value = "Electronique"
$ws.Range("F5").Value = "semestre 2"
$ws.Range("G5").Value = 2019
$ws.Range("H5").Value = "Absent(e)"
